# Auto-generated edit script applying numeric corrections to the
# per-item profit-tracking sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Each block updates price/profit columns (H,I,J,K,L,M,N) for one row.
$wb = $excel.ActiveWorkbook

# ALC row 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45000
$ws.Range("J3").Value = 45000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45228

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11750.25
$ws.Range("J40").Value = 11666.667
$ws.Range("L40").Value = 11666.667
$ws.Range("N40").Value = -12016.667

# ALC row 57
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 53799.8
$ws.Range("J57").Value = 51999.75
$ws.Range("L57").Value = 155999.25
$ws.Range("N57").Value = -156997.25

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8000
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25748

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 8000
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -80736

# ALC row 102
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 45000
$ws.Range("J102").Value = 45000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -51490

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1942.1875
$ws.Range("I112").Value = 1400
$ws.Range("K112").Value = 4200
$ws.Range("M112").Value = -3092

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5237.2627
$ws.Range("J138").Value = 6025.3887
$ws.Range("L138").Value = 18076.1661
$ws.Range("N138").Value = -28356.1661

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11240.7
$ws.Range("I32").Value = 5549.5713
$ws.Range("K32").Value = 5549.5713
$ws.Range("M32").Value = -5262.5713

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3620.1018
$ws.Range("I61").Value = 3191.0728
$ws.Range("K61").Value = 3191.0728
$ws.Range("M61").Value = -2979.0728

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1619.9706
$ws.Range("I102").Value = 1660.1072
$ws.Range("J102").Value = 1432.6666
$ws.Range("K102").Value = 1660.1072
$ws.Range("L102").Value = 1432.6666
$ws.Range("M102").Value = -38.10719999999992
$ws.Range("N102").Value = -4676.6666

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 180341.97
$ws.Range("I110").Value = 201721.05
$ws.Range("K110").Value = 201721.05
$ws.Range("M110").Value = -199676.05

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2191.0833
$ws.Range("I122").Value = 1667.1111
$ws.Range("K122").Value = 5001.3333
$ws.Range("M122").Value = -2551.3333

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3838.0527
$ws.Range("I132").Value = 2884.8071
$ws.Range("K132").Value = 8654.4213
$ws.Range("M132").Value = -6124.4213

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3620.1018
$ws.Range("I136").Value = 3191.0728
$ws.Range("K136").Value = 9573.2184
$ws.Range("M136").Value = -7023.2184

# BSM row 12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 5948
$ws.Range("I12").Value = 5948
$ws.Range("K12").Value = 5948
$ws.Range("M12").Value = -5780

# BSM row 17
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = 0

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1972.1578
$ws.Range("J20").Value = 1881.1666
$ws.Range("L20").Value = 1881.1666
$ws.Range("N20").Value = -2375.1666

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2824.6667
$ws.Range("I99").Value = 3237
$ws.Range("K99").Value = 3237
$ws.Range("M99").Value = -1739

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 72599.71000000001
$ws.Range("I31").Value = 1261.2307
$ws.Range("K31").Value = 1261.2307
$ws.Range("M31").Value = -966.2307000000001

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 72599.71000000001
$ws.Range("I34").Value = 1261.2307
$ws.Range("K34").Value = 1261.2307
$ws.Range("M34").Value = -1059.2307

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2006.7894
$ws.Range("I58").Value = 1854.75
$ws.Range("K58").Value = 1854.75
$ws.Range("M58").Value = -1651.75

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4649.6665
$ws.Range("I62").Value = 4474.75
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 4474.75
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -3850.75
$ws.Range("N62").Value = -6247.5

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4649.6665
$ws.Range("I65").Value = 4474.75
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 22373.75
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -19253.75
$ws.Range("N65").Value = -31237.5

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1132
$ws.Range("J105").Value = 999
$ws.Range("L105").Value = 999
$ws.Range("N105").Value = -4493

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2006.7894
$ws.Range("I136").Value = 1854.75
$ws.Range("K136").Value = 5564.25
$ws.Range("M136").Value = -3014.25

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1770.5
$ws.Range("I92").Value = 1229.6
$ws.Range("J92").Value = 1978.5385
$ws.Range("K92").Value = 3688.8
$ws.Range("L92").Value = 5935.6155
$ws.Range("M92").Value = -2440.8
$ws.Range("N92").Value = -8431.6155

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 13443.5
$ws.Range("J107").Value = 25779.5
$ws.Range("L107").Value = 77338.5
$ws.Range("N107").Value = -81178.5

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5853.364
$ws.Range("I139").Value = 4399.6
$ws.Range("K139").Value = 13198.8
$ws.Range("M139").Value = -8058.800000000001

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4331.5
$ws.Range("I140").Value = 2798.6
$ws.Range("K140").Value = 8395.799999999999
$ws.Range("M140").Value = -3215.799999999999

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 7113.9375
$ws.Range("J141").Value = 8999.4
$ws.Range("L141").Value = 26998.2
$ws.Range("N141").Value = -37358.2

# GSM row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10001.5
$ws.Range("J4").Value = 10001.5
$ws.Range("L4").Value = 10001.5
$ws.Range("N4").Value = -10225.5

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 71438104
$ws.Range("I70").Value = 7147.5
$ws.Range("K70").Value = 7147.5
$ws.Range("M70").Value = -6877.5

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 71438104
$ws.Range("I73").Value = 7147.5
$ws.Range("K73").Value = 7147.5
$ws.Range("M73").Value = -6211.5

# GSM row 82
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# GSM row 85
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 30345.324
$ws.Range("I132").Value = 5737
$ws.Range("J132").Value = 81454.92
$ws.Range("K132").Value = 17211
$ws.Range("L132").Value = 244364.76
$ws.Range("M132").Value = -14681
$ws.Range("N132").Value = -249424.76

# LTW row 8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H8").Value = 77500
$ws.Range("J8").Value = 77500
$ws.Range("L8").Value = 77500
$ws.Range("N8").Value = -77780

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7656.879
$ws.Range("I61").Value = 6974.25
$ws.Range("J61").Value = 11479.6
$ws.Range("K61").Value = 6974.25
$ws.Range("L61").Value = 11479.6
$ws.Range("M61").Value = -6772.25
$ws.Range("N61").Value = -11883.6

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2049.5
$ws.Range("I100").Value = 1499.5
$ws.Range("K100").Value = 1499.5
$ws.Range("M100").Value = -958.5

# LTW row 110
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 52665
$ws.Range("J110").Value = 52665
$ws.Range("L110").Value = 52665
$ws.Range("N110").Value = -60845

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7656.879
$ws.Range("I113").Value = 6974.25
$ws.Range("J113").Value = 11479.6
$ws.Range("K113").Value = 6974.25
$ws.Range("L113").Value = 11479.6
$ws.Range("M113").Value = -4804.25
$ws.Range("N113").Value = -15819.6

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4062.658
$ws.Range("I132").Value = 4098.1934
$ws.Range("K132").Value = 12294.5802
$ws.Range("M132").Value = -9764.5802

# WVR row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 10118
$ws.Range("I75").Value = 10118
$ws.Range("K75").Value = 10118
$ws.Range("M75").Value = -9182

# WVR row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 10118
$ws.Range("I78").Value = 10118
$ws.Range("K78").Value = 30354
$ws.Range("M78").Value = -25674

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 44756.22
$ws.Range("J96").Value = 1466.5
$ws.Range("L96").Value = 1466.5
$ws.Range("N96").Value = -4212.5

# WVR row 101
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 50602
$ws.Range("J101").Value = 50602
$ws.Range("L101").Value = 50602
$ws.Range("N101").Value = -57092

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41400.19
$ws.Range("I132").Value = 2619.5557
$ws.Range("K132").Value = 7858.6671
$ws.Range("M132").Value = -5328.6671

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 292289.66
$ws.Range("I136").Value = 309913.1
$ws.Range("J136").Value = 225712.22
$ws.Range("K136").Value = 929739.2999999999
$ws.Range("L136").Value = 677136.66
$ws.Range("M136").Value = -927189.2999999999
$ws.Range("N136").Value = -682236.66
